$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78. This pushes the existing rows 78-127
# down to 79-128, preserving all of their data (matches the diff, which
# shows every row from 78 to 127 shifting its values down by one row and
# a brand-new row 128 appearing with what used to be row 127's values).
$ws.Rows(78).Insert()

# Populate the newly inserted row 78 with the new record's data. The
# "constant" columns (A, B, C, E, F, G, H, I, R) are identical for every
# row in this sheet, so copy them from the row directly below (row 79,
# which now holds the data that used to live in row 78).
$ws.Range("A78").Value = $ws.Range("A79").Value()
$ws.Range("B78").Value = $ws.Range("B79").Value()
$ws.Range("C78").Value = $ws.Range("C79").Value()
$ws.Range("D78").Value = 44719
$ws.Range("E78").Value = $ws.Range("E79").Value()
$ws.Range("F78").Value = $ws.Range("F79").Value()
$ws.Range("G78").Value = $ws.Range("G79").Value()
$ws.Range("H78").Value = $ws.Range("H79").Value()
$ws.Range("I78").Value = $ws.Range("I79").Value()
$ws.Range("J78").Value = 100
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 10000
$ws.Range("M78").Value = 9500
$ws.Range("N78").Value = '$/caja 50 unidades'
$ws.Range("O78").Value = 'Región de Arica y Parinacota'
$ws.Range("P78").Value = 190
$ws.Range("Q78").Value = 50
$ws.Range("R78").Value = $ws.Range("R79").Value()

# Match the date-formatted number format used by the other rows' Fecha column.
$ws.Range("D78").NumberFormat = $ws.Range("D79").NumberFormat()
